# Rw Daily Usage - 6:55PM
# Updates the Dashboard "spend tracker" grid (rows 4-11) with the day's
# latest tallies, logs the new timestamp for each settled row, appends a
# new Purchase entry (row 17 / sheet row 26), and records the active cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Dashboard sheet
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

# Reserved cash split changes -> C2 (=B3+B4) recalculates automatically.
$dash.Range("B4").Value = 125

# Row 4
$dash.Range("E4").Value = 20
$dash.Range("F4").ClearContents()
$dash.Range("L4").Value = "ABP"
$dash.Range("N4").Value = 0.78472222222222221
$dash.Range("P4").ClearContents()

# Row 5
$dash.Range("E5").Value = 40
$dash.Range("F5").ClearContents()
$dash.Range("N5").Value = 0.78472222222222221
$dash.Range("P5").ClearContents()

# Row 6
$dash.Range("F6").ClearContents()
$dash.Range("L6").Value = "Rbharat"
$dash.Range("N6").Value = 0.78472222222222221
$dash.Range("P6").ClearContents()

# Row 7
$dash.Range("E7").Value = 10
$dash.Range("F7").ClearContents()
$dash.Range("N7").Value = 0.78472222222222221
$dash.Range("O7").Value = "g0"
$dash.Range("P7").ClearContents()

# Row 8
$dash.Range("E8").Value = 50
$dash.Range("N8").Value = 0.78472222222222221
$dash.Range("P8").ClearContents()

# Row 9 (B9 = 7000-C2 recalculates automatically)
$dash.Range("E9").ClearContents()
$dash.Range("F9").Value = 5
$dash.Range("N9").Value = 0.78472222222222221
$dash.Range("P9").ClearContents()

# Row 10
$dash.Range("E10").ClearContents()
$dash.Range("F10").Value = 30
$dash.Range("I10").Value = "MatriX"
$dash.Range("N10").Value = 0.78472222222222221
$dash.Range("P10").ClearContents()

# Row 11
$dash.Range("E11").Value = 10
$dash.Range("F11").ClearContents()
$dash.Range("N11").Value = 0.78472222222222221
$dash.Range("O11").Value = "pp"
$dash.Range("P11").ClearContents()

# G12 (=SUM(E4:F11)), B13 (=B18+Purchase!O2) recalc automatically.

# ---------------------------------------------------------------------
# Purchase sheet - add today's entry (row 17 / sheet row 26)
# ---------------------------------------------------------------------
$purchase = $wb.Worksheets.Item("Purchase")

$purchase.Range("B26").Value = 17
$purchase.Range("E26").Value = 35
$purchase.Range("F26").Value = 50
$purchase.Range("G26").Value = 40
$purchase.Range("H26").Value = 10
$purchase.Range("I26").Value = 20
$purchase.Range("J26").Value = 20

# Match the row-26 formatting to the row above it (row 25).
$purchase.Range("B25:M25").Copy()
$purchase.Range("B26:M26").PasteSpecial(-4122)
$purchase.Range("B26").Value = 17
$purchase.Range("E26").Value = 35
$purchase.Range("F26").Value = 50
$purchase.Range("G26").Value = 40
$purchase.Range("H26").Value = 10
$purchase.Range("I26").Value = 20
$purchase.Range("J26").Value = 20
$purchase.Range("C26").ClearContents()
$purchase.Range("D26").ClearContents()
$purchase.Range("K26").ClearContents()
$purchase.Range("L26").ClearContents()
$purchase.Range("M26").ClearContents()

# Extend the running total to include the new row.
$purchase.Range("O2").Formula = "=SUM(E2:J26)"

# ---------------------------------------------------------------------
# Leave the selection on the Dashboard sheet where the user last clicked.
# ---------------------------------------------------------------------
$dash.Activate()
$dash.Range("M16").Select()
